# TeamLab-presentation.pptx - "new folder for zip-submission"
#
# This edit stamps four slides (the SmartArt/process slides 6-9) with a
# custom "TIMING" tag used by the presenter's timer tooling. Adding a
# tag creates ppt/tags/tagN.xml + the r:id relationship + the
# <p:custDataLst> element on each slide automatically.
#
# (The deck's Notes Master also re-caches its auto date field from
# 06.06.2021 to 07.06.2021 when PowerPoint re-saves it a day later, but
# that field lives on the Notes Master, which this host's object model
# does not expose a safe/working write path for -- Notes Master shape
# writes here misroute onto the Slide Master, so that part of the edit
# is intentionally left alone rather than risk corrupting the deck.)

$p = $ppt.ActivePresentation

# TIMING tags for slides 6-9 (order matters: tag1..tag4 are minted in
# the order the Tags.Add calls run).
$p.Slides.Item(6).Tags.Add("TIMING", "|7.8|2|2.5|2.3|6")
$p.Slides.Item(7).Tags.Add("TIMING", "|8.7|7.6")
$p.Slides.Item(8).Tags.Add("TIMING", "|5.1|12.6|16.5")
$p.Slides.Item(9).Tags.Add("TIMING", "|6|10.7")
